$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "G2"; Value = 189.0573523333333 },
    @{ Cell = "H2"; Value = 567.172057 },
    @{ Cell = "I2"; Value = 0.1182556374491171 },
    @{ Cell = "J2"; Value = 0.1182556374491171 },
    @{ Cell = "M2"; Value = 11.42350833333333 },
    @{ Cell = "N2"; Value = 34.270525 },
    @{ Cell = "O2"; Value = 0.2763488058447062 },
    @{ Cell = "P2"; Value = 0.2763488058447063 },
    @{ Cell = "Q2"; Value = 2159.698239857769 },
    @{ Cell = "R2"; Value = 19437.28415871993 },
    @{ Cell = "S2"; Value = 0.03267980419346803 },
    @{ Cell = "T2"; Value = 0.03267980419346803 },
    @{ Cell = "G3"; Value = 189.0573523333333 },
    @{ Cell = "H3"; Value = 567.172057 },
    @{ Cell = "I3"; Value = 0.1182556374491171 },
    @{ Cell = "J3"; Value = 0.1182556374491171 },
    @{ Cell = "O3"; Value = 0.04096986897477309 },
    @{ Cell = "P3"; Value = 0.04096986897477309 },
    @{ Cell = "Q3"; Value = 320.1843179367437 },
    @{ Cell = "R3"; Value = 2881.658861430693 },
    @{ Cell = "S3"; Value = 0.004844917971818597 },
    @{ Cell = "T3"; Value = 0.004844917971818597 },
    @{ Cell = "G4"; Value = 189.0573523333333 },
    @{ Cell = "H4"; Value = 567.172057 },
    @{ Cell = "I4"; Value = 0.1182556374491171 },
    @{ Cell = "J4"; Value = 0.1182556374491171 },
    @{ Cell = "M4"; Value = 8.972155666666666 },
    @{ Cell = "N4"; Value = 26.916467 },
    @{ Cell = "O4"; Value = 0.2170475507162041 },
    @{ Cell = "P4"; Value = 0.2170475507162042 },
    @{ Cell = "Q4"; Value = 1696.251995062513 },
    @{ Cell = "R4"; Value = 15266.26795556262 },
    @{ Cell = "S4"; Value = 0.02566709646671429 },
    @{ Cell = "T4"; Value = 0.02566709646671429 },
    @{ Cell = "G5"; Value = 189.0573523333333 },
    @{ Cell = "H5"; Value = 567.172057 },
    @{ Cell = "I5"; Value = 0.1182556374491171 },
    @{ Cell = "J5"; Value = 0.1182556374491171 },
    @{ Cell = "M5"; Value = 19.24803433333333 },
    @{ Cell = "N5"; Value = 57.744103 },
    @{ Cell = "O5"; Value = 0.4656337744643164 },
    @{ Cell = "P5"; Value = 0.4656337744643165 },
    @{ Cell = "Q5"; Value = 3638.982408681097 },
    @{ Cell = "R5"; Value = 32750.84167812987 },
    @{ Cell = "S5"; Value = 0.05506381881711616 },
    @{ Cell = "T5"; Value = 0.05506381881711616 },
    @{ Cell = "G6"; Value = 930.1503093333332 },
    @{ Cell = "I6"; Value = 0.5818103152093762 },
    @{ Cell = "J6"; Value = 0.5818103152093762 },
    @{ Cell = "M6"; Value = 11.42350833333333 },
    @{ Cell = "N6"; Value = 34.270525 },
    @{ Cell = "O6"; Value = 0.2763488058447062 },
    @{ Cell = "P6"; Value = 0.2763488058447063 },
    @{ Cell = "Q6"; Value = 10625.57980992191 },
    @{ Cell = "R6"; Value = 95630.21828929719 },
    @{ Cell = "S6"; Value = 0.1607825858362432 },
    @{ Cell = "T6"; Value = 0.1607825858362433 },
    @{ Cell = "G7"; Value = 930.1503093333332 },
    @{ Cell = "I7"; Value = 0.5818103152093762 },
    @{ Cell = "J7"; Value = 0.5818103152093762 },
    @{ Cell = "O7"; Value = 0.04096986897477309 },
    @{ Cell = "P7"; Value = 0.04096986897477309 },
    @{ Cell = "S7"; Value = 0.02383669238229957 },
    @{ Cell = "T7"; Value = 0.02383669238229957 },
    @{ Cell = "G8"; Value = 930.1503093333332 },
    @{ Cell = "I8"; Value = 0.5818103152093762 },
    @{ Cell = "J8"; Value = 0.5818103152093762 },
    @{ Cell = "M8"; Value = 8.972155666666666 },
    @{ Cell = "N8"; Value = 26.916467 },
    @{ Cell = "O8"; Value = 0.2170475507162041 },
    @{ Cell = "P8"; Value = 0.2170475507162042 },
    @{ Cell = "Q8"; Value = 8345.453368736818 },
    @{ Cell = "R8"; Value = 75109.08031863136 },
    @{ Cell = "S8"; Value = 0.1262805038976178 },
    @{ Cell = "T8"; Value = 0.1262805038976178 },
    @{ Cell = "G9"; Value = 930.1503093333332 },
    @{ Cell = "I9"; Value = 0.5818103152093762 },
    @{ Cell = "J9"; Value = 0.5818103152093762 },
    @{ Cell = "M9"; Value = 19.24803433333333 },
    @{ Cell = "N9"; Value = 57.744103 },
    @{ Cell = "O9"; Value = 0.4656337744643164 },
    @{ Cell = "P9"; Value = 0.4656337744643165 },
    @{ Cell = "Q9"; Value = 17903.56508920862 },
    @{ Cell = "R9"; Value = 161132.0858028776 },
    @{ Cell = "S9"; Value = 0.2709105330932155 },
    @{ Cell = "T9"; Value = 0.2709105330932156 },
    @{ Cell = "G10"; Value = 420.6651306666666 },
    @{ Cell = "H10"; Value = 1261.995392 },
    @{ Cell = "I10"; Value = 0.2631266256807295 },
    @{ Cell = "J10"; Value = 0.2631266256807295 },
    @{ Cell = "M10"; Value = 11.42350833333333 },
    @{ Cell = "N10"; Value = 34.270525 },
    @{ Cell = "O10"; Value = 0.2763488058447062 },
    @{ Cell = "P10"; Value = 0.2763488058447063 },
    @{ Cell = "Q10"; Value = 4805.471625713421 },
    @{ Cell = "R10"; Value = 43249.24463142079 },
    @{ Cell = "S10"; Value = 0.0727147287928166 },
    @{ Cell = "T10"; Value = 0.07271472879281662 },
    @{ Cell = "G11"; Value = 420.6651306666666 },
    @{ Cell = "H11"; Value = 1261.995392 },
    @{ Cell = "I11"; Value = 0.2631266256807295 },
    @{ Cell = "J11"; Value = 0.2631266256807295 },
    @{ Cell = "O11"; Value = 0.04096986897477309 },
    @{ Cell = "P11"; Value = 0.04096986897477309 },
    @{ Cell = "Q11"; Value = 712.4313139898452 },
    @{ Cell = "R11"; Value = 6411.881825908607 },
    @{ Cell = "S11"; Value = 0.01078026337791365 },
    @{ Cell = "T11"; Value = 0.01078026337791365 },
    @{ Cell = "G12"; Value = 420.6651306666666 },
    @{ Cell = "H12"; Value = 1261.995392 },
    @{ Cell = "I12"; Value = 0.2631266256807295 },
    @{ Cell = "J12"; Value = 0.2631266256807295 },
    @{ Cell = "M12"; Value = 8.972155666666666 },
    @{ Cell = "N12"; Value = 26.916467 },
    @{ Cell = "O12"; Value = 0.2170475507162041 },
    @{ Cell = "P12"; Value = 0.2170475507162042 },
    @{ Cell = "Q12"; Value = 3774.273035880006 },
    @{ Cell = "R12"; Value = 33968.45732292005 },
    @{ Cell = "S12"; Value = 0.05711098963222179 },
    @{ Cell = "T12"; Value = 0.05711098963222181 },
    @{ Cell = "G13"; Value = 420.6651306666666 },
    @{ Cell = "H13"; Value = 1261.995392 },
    @{ Cell = "I13"; Value = 0.2631266256807295 },
    @{ Cell = "J13"; Value = 0.2631266256807295 },
    @{ Cell = "M13"; Value = 19.24803433333333 },
    @{ Cell = "N13"; Value = 57.744103 },
    @{ Cell = "O13"; Value = 0.4656337744643164 },
    @{ Cell = "P13"; Value = 0.4656337744643165 },
    @{ Cell = "Q13"; Value = 8096.976877908151 },
    @{ Cell = "R13"; Value = 72872.79190117336 },
    @{ Cell = "S13"; Value = 0.1225206438777774 },
    @{ Cell = "T13"; Value = 0.1225206438777774 },
    @{ Cell = "G14"; Value = 58.84466766666667 },
    @{ Cell = "H14"; Value = 176.534003 },
    @{ Cell = "I14"; Value = 0.03680742166077718 },
    @{ Cell = "J14"; Value = 0.03680742166077718 },
    @{ Cell = "M14"; Value = 11.42350833333333 },
    @{ Cell = "N14"; Value = 34.270525 },
    @{ Cell = "O14"; Value = 0.2763488058447062 },
    @{ Cell = "P14"; Value = 0.2763488058447063 },
    @{ Cell = "Q14"; Value = 672.2125514623973 },
    @{ Cell = "R14"; Value = 6049.912963161575 },
    @{ Cell = "S14"; Value = 0.01017168702217835 },
    @{ Cell = "T14"; Value = 0.01017168702217835 },
    @{ Cell = "G15"; Value = 58.84466766666667 },
    @{ Cell = "H15"; Value = 176.534003 },
    @{ Cell = "I15"; Value = 0.03680742166077718 },
    @{ Cell = "J15"; Value = 0.03680742166077718 },
    @{ Cell = "O15"; Value = 0.04096986897477309 },
    @{ Cell = "P15"; Value = 0.04096986897477309 },
    @{ Cell = "Q15"; Value = 99.65832880091635 },
    @{ Cell = "R15"; Value = 896.9249592082471 },
    @{ Cell = "S15"; Value = 0.001507995242741266 },
    @{ Cell = "T15"; Value = 0.001507995242741266 },
    @{ Cell = "G16"; Value = 58.84466766666667 },
    @{ Cell = "H16"; Value = 176.534003 },
    @{ Cell = "I16"; Value = 0.03680742166077718 },
    @{ Cell = "J16"; Value = 0.03680742166077718 },
    @{ Cell = "M16"; Value = 8.972155666666666 },
    @{ Cell = "N16"; Value = 26.916467 },
    @{ Cell = "O16"; Value = 0.2170475507162041 },
    @{ Cell = "P16"; Value = 0.2170475507162042 },
    @{ Cell = "Q16"; Value = 527.9635184586001 },
    @{ Cell = "R16"; Value = 4751.671666127401 },
    @{ Cell = "S16"; Value = 0.007988960719650245 },
    @{ Cell = "T16"; Value = 0.007988960719650247 },
    @{ Cell = "G17"; Value = 58.84466766666667 },
    @{ Cell = "H17"; Value = 176.534003 },
    @{ Cell = "I17"; Value = 0.03680742166077718 },
    @{ Cell = "J17"; Value = 0.03680742166077718 },
    @{ Cell = "M17"; Value = 19.24803433333333 },
    @{ Cell = "N17"; Value = 57.744103 },
    @{ Cell = "O17"; Value = 0.4656337744643164 },
    @{ Cell = "P17"; Value = 0.4656337744643165 },
    @{ Cell = "Q17"; Value = 1132.64418358159 },
    @{ Cell = "R17"; Value = 10193.79765223431 },
    @{ Cell = "S17"; Value = 0.01713877867620732 },
    @{ Cell = "T17"; Value = 0.01713877867620732 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Output "Applied $($updates.Count) cell updates"
